$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: date and volume update ---
$ws.Range("D18").Value = 44553
$ws.Range("M18").Value = 250

# --- Row 19: date and volume update ---
$ws.Range("D19").Value = 44553
$ws.Range("M19").Value = 250

# --- Row 20: date and volume update ---
$ws.Range("D20").Value = 44553
$ws.Range("M20").Value = 250

# --- Row 21: becomes "Especial" record (previously "Primera") ---
$ws.Range("D21").Value = 44551
$ws.Range("L21").Value = "Especial"
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 20000
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 20000
$ws.Range("Q21").Value = "$/caja 18 kilos"
$ws.Range("S21").Value = 1111
$ws.Range("T21").Value = 18

# --- Row 22: becomes "Primera" record (previously "Segunda") ---
$ws.Range("D22").Value = 44551
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = 18000
$ws.Range("O22").Value = 18000
$ws.Range("P22").Value = 18000
$ws.Range("Q22").Value = "$/caja 18 kilos"
$ws.Range("S22").Value = 1000
$ws.Range("T22").Value = 18

# --- Row 23: stays "Segunda", values updated to the 18kg figures ---
$ws.Range("D23").Value = 44551
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 16000
$ws.Range("P23").Value = 16000
$ws.Range("Q23").Value = "$/caja 18 kilos"
$ws.Range("S23").Value = 889
$ws.Range("T23").Value = 18

# --- Row 24: becomes "Primera" record with the old row-21 figures ---
$ws.Range("D24").Value = 44187
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 350
$ws.Range("Q24").Value = "$/caja 15 kilos"
$ws.Range("R24").Value = "Región Metropolitana"
$ws.Range("S24").Value = 1067
$ws.Range("T24").Value = 15

# --- New row 25 (re-insertion of the old row-22 "Segunda" record) ---
$ws.Range("A25").Value = 4
$ws.Range("B25").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C25").Value = "Los Lagos"
$ws.Range("D25").Value = 44187
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100103
$ws.Range("H25").Value = "Frutos de hueso (carozo)"
$ws.Range("I25").Value = 100103003
$ws.Range("J25").Value = "Damasco"
$ws.Range("K25").Value = "Castle Brite"
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 300
$ws.Range("N25").Value = 13000
$ws.Range("O25").Value = 13000
$ws.Range("P25").Value = 13000
$ws.Range("Q25").Value = "$/caja 15 kilos"
$ws.Range("R25").Value = "Región Metropolitana"
$ws.Range("S25").Value = 867
$ws.Range("T25").Value = 15

# --- New row 26 (re-insertion of the old row-23 "Segunda" record) ---
$ws.Range("A26").Value = 4
$ws.Range("B26").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C26").Value = "Los Lagos"
$ws.Range("D26").Value = 44194
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100103
$ws.Range("H26").Value = "Frutos de hueso (carozo)"
$ws.Range("I26").Value = 100103003
$ws.Range("J26").Value = "Damasco"
$ws.Range("K26").Value = "Castle Brite"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 300
$ws.Range("N26").Value = 15000
$ws.Range("O26").Value = 16000
$ws.Range("P26").Value = 15500
$ws.Range("Q26").Value = "$/caja 15 kilos"
$ws.Range("R26").Value = "Región Metropolitana"
$ws.Range("S26").Value = 1033
$ws.Range("T26").Value = 15

# --- New row 27 (re-insertion of the old row-24 "Segunda" record) ---
$ws.Range("A27").Value = 4
$ws.Range("B27").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C27").Value = "Los Lagos"
$ws.Range("D27").Value = 44540
$ws.Range("E27").Value = 10
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100103
$ws.Range("H27").Value = "Frutos de hueso (carozo)"
$ws.Range("I27").Value = 100103003
$ws.Range("J27").Value = "Damasco"
$ws.Range("K27").Value = "Castle Brite"
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 600
$ws.Range("N27").Value = 16000
$ws.Range("O27").Value = 16000
$ws.Range("P27").Value = 16000
$ws.Range("Q27").Value = "$/caja 18 kilos"
$ws.Range("R27").Value = "Región del Maule"
$ws.Range("S27").Value = 889
$ws.Range("T27").Value = 18

# Keep the date-formatted column consistent with the rest of column D.
$ws.Range("D25:D27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
